# Auto-generated edit script applying the cryptos.xlsx crypto-price refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D4 ("TetherUSD" price, "1.00") is never touched by this edit and keeps the
# sheets default (unstyled) cell format. We use it as a style template so that
# forcing numeric-looking price strings to stay text (via NumberFormat "@") does
# not leave a stray custom style behind once the value has been written.
$defaultStyle = $ws.Cells.Item(4, 4).Style

$ws.Cells.Item(2, 4).Value = "66.178.75"
$ws.Cells.Item(2, 5).Value = "  -4.73%  "

$ws.Cells.Item(3, 4).Value = "3.277.15"
$ws.Cells.Item(3, 5).Value = "  -5.85%  "

$ws.Cells.Item(4, 5).Value = "  +0.09%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "558.87"
$ws.Cells.Item(5, 4).Style = $defaultStyle

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "186.28"
$ws.Cells.Item(6, 4).Style = $defaultStyle
$ws.Cells.Item(6, 5).Value = "  -3.54%  "

$ws.Cells.Item(7, 5).Value = "  +0.01%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.591"
$ws.Cells.Item(8, 4).Style = $defaultStyle
$ws.Cells.Item(8, 5).Value = "  -3.01%  "

$ws.Cells.Item(9, 4).Value = "3.273.23"
$ws.Cells.Item(9, 5).Value = "  -5.66%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.187"
$ws.Cells.Item(10, 4).Style = $defaultStyle
$ws.Cells.Item(10, 5).Value = "  -9.04%  "

$ws.Cells.Item(11, 5).Value = "  -4.98%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "47.52"
$ws.Cells.Item(12, 4).Style = $defaultStyle
$ws.Cells.Item(12, 5).Value = "  -7.43%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.0000266"
$ws.Cells.Item(13, 4).Style = $defaultStyle
$ws.Cells.Item(13, 5).Value = "  -6.71%  "

$ws.Cells.Item(14, 2).Value = "BitcoinCash"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "634.48"
$ws.Cells.Item(14, 4).Style = $defaultStyle
$ws.Cells.Item(14, 5).Value = "  -1.28%  "

$ws.Cells.Item(15, 2).Value = "Polkadot"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "8.60"
$ws.Cells.Item(15, 4).Style = $defaultStyle
$ws.Cells.Item(15, 5).Value = "  -5.89%  "

$ws.Cells.Item(16, 4).Value = "3.803.76"
$ws.Cells.Item(16, 5).Value = "  -5.93%  "

$ws.Cells.Item(17, 4).Value = "66.083.21"
$ws.Cells.Item(17, 5).Value = "  -4.51%  "

$ws.Cells.Item(18, 5).Value = "  -1.49%  "

$ws.Cells.Item(19, 5).Value = "  -3.58%  "

$ws.Cells.Item(20, 4).Value = "3.280.65"
$ws.Cells.Item(20, 5).Value = "  -5.51%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "11.36"
$ws.Cells.Item(21, 4).Style = $defaultStyle
$ws.Cells.Item(21, 5).Value = "  -8.02%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.905"
$ws.Cells.Item(22, 4).Style = $defaultStyle
$ws.Cells.Item(22, 5).Value = "  -4.28%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "18.37"
$ws.Cells.Item(23, 4).Style = $defaultStyle
$ws.Cells.Item(23, 5).Value = "  +3.27%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "107.56"
$ws.Cells.Item(24, 4).Style = $defaultStyle
$ws.Cells.Item(24, 5).Value = "  +8.56%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "4.91"
$ws.Cells.Item(25, 4).Style = $defaultStyle
$ws.Cells.Item(25, 5).Value = "  -7.27%  "

$ws.Cells.Item(26, 5).Value = "  -7.23%  "

$ws.Cells.Item(27, 5).Value = "  -7.13%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "9.61"
$ws.Cells.Item(28, 4).Style = $defaultStyle
$ws.Cells.Item(28, 5).Value = "  -3.42%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "8.71"
$ws.Cells.Item(29, 4).Style = $defaultStyle
$ws.Cells.Item(29, 5).Value = "  -6.50%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "30.33"
$ws.Cells.Item(30, 4).Style = $defaultStyle
$ws.Cells.Item(30, 5).Value = "  -6.67%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "4.08"
$ws.Cells.Item(31, 4).Style = $defaultStyle
$ws.Cells.Item(31, 5).Value = "  -5.01%  "

$ws.Cells.Item(32, 5).Value = "  -6.92%  "

$ws.Cells.Item(33, 5).Value = "  -4.99%  "

$ws.Cells.Item(34, 5).Value = "  -3.97%  "

$ws.Cells.Item(35, 2).Value = "OKB"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "57.68"
$ws.Cells.Item(35, 4).Style = $defaultStyle
$ws.Cells.Item(35, 5).Value = "  -5.33%  "

$ws.Cells.Item(36, 2).Value = "Bittensor"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "530.07"
$ws.Cells.Item(36, 4).Style = $defaultStyle
$ws.Cells.Item(36, 5).Value = "  +1.79%  "

$ws.Cells.Item(37, 2).Value = "Maker"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(37, 4).Value = "3.734.11"
$ws.Cells.Item(37, 5).Value = "  -0.14%  "

$ws.Cells.Item(38, 5).Value = "  -0.07%  "

$ws.Cells.Item(39, 5).Value = "  -3.94%  "

$ws.Cells.Item(40, 4).Value = "0.0₃0725"
$ws.Cells.Item(40, 5).Value = "  -8.51%  "

$ws.Cells.Item(41, 2).Value = "Fetch.AI"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "2.75"
$ws.Cells.Item(41, 4).Style = $defaultStyle
$ws.Cells.Item(41, 5).Value = "  -6.87%  "

$ws.Cells.Item(42, 2).Value = "CoreDAO"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "3.43"
$ws.Cells.Item(42, 4).Style = $defaultStyle
$ws.Cells.Item(42, 5).Value = "  -2.23%  "

$ws.Cells.Item(43, 2).Value = "Kaspa"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.129"
$ws.Cells.Item(43, 4).Style = $defaultStyle
$ws.Cells.Item(43, 5).Value = "  -3.48%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "32.91"
$ws.Cells.Item(44, 4).Style = $defaultStyle
$ws.Cells.Item(44, 5).Value = "  -4.11%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.338"
$ws.Cells.Item(45, 4).Style = $defaultStyle
$ws.Cells.Item(45, 5).Value = "  -9.00%  "

$ws.Cells.Item(46, 5).Value = "  -2.03%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.0415"
$ws.Cells.Item(47, 4).Style = $defaultStyle
$ws.Cells.Item(47, 5).Value = "  -6.31%  "

$ws.Cells.Item(48, 2).Value = "Stellar"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.129"
$ws.Cells.Item(48, 4).Style = $defaultStyle
$ws.Cells.Item(48, 5).Value = "  -3.64%  "

$ws.Cells.Item(49, 2).Value = "ThetaToken"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "2.61"
$ws.Cells.Item(49, 4).Style = $defaultStyle
$ws.Cells.Item(49, 5).Value = "  -7.50%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.999"
$ws.Cells.Item(50, 4).Style = $defaultStyle
$ws.Cells.Item(50, 5).Value = "  +0.02%  "

$ws.Cells.Item(51, 5).Value = "  +1.99%  "
